# Updated input forms with dummy template values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows -------------------------------------------------
# Row 2: server URL / project name / PAT -> dummy template values
$ws.Range("A2").Value = "http://128.0.0.1/TestCollection"
$ws.Range("B2").Value = "project1"
$ws.Range("C2").Value = "adad87adad8ds4449m434344mmnbnbb43434"

# Row 3: server URL / project name / PAT -> dummy template values
$ws.Range("A3").Value = "http://128.0.0.1/TestCollection"
$ws.Range("B3").Value = "project2"
$ws.Range("C3").Value = "adad87adad8ds4449m434344mmnbnbb43434"

# Remove the now-unused extra rows (devserver / testproj)
$ws.Rows("4:5").Delete()

# --- Fix up hyperlinks --------------------------------------------------
# Deleting via a range wipes the whole worksheet hyperlink collection in
# this engine, so do it once, then re-add the two that should remain,
# preserving each cell's original style afterwards.
$styleA2 = $ws.Range("A2").Style
$styleA3 = $ws.Range("A3").Style

$ws.Range("A2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "http://128.0.0.1/TestCollection")
$ws.Hyperlinks.Add($ws.Range("A3"), "http://128.0.0.1/TestCollection")

$ws.Range("A2").Style = $styleA2
$ws.Range("A3").Style = $styleA3

# Clear the stale selection left over from the deleted rows
$ws.Range("A1").Select()
